$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column N (year 2023) ------------------------------------
# Row 3 (thin border separator row) - extend to N3
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)

# Row 4 (year headers) - N4 = 2023
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2023

# Row 5 (Small enterprises / Чакан ишканалар data) - N5
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 2.5449890821474286

# Row 6 (Medium-sized enterprises / Орто ишканалар data) - N6
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 1.4569686017619159

# --- Fix the Kyrgyz title in A1: drop the trailing newline ------------
$ws.Range("A1").Value = "8.3.1.2 Экономикадагы иш менен камсыз болгон бардык калктын чакан жана орто ишканаларда иштегендердин үлүшү"

# --- Row height tweaks --------------------------------------------------
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 17.25
$ws.Rows.Item(6).RowHeight = 17.25
